$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (scott) - updated age/payment
$ws.Range("B2").Value = 24
$ws.Range("C2").Value = 406

# Row 11 (morgan) - updated age/payment
$ws.Range("B11").Value = 12
$ws.Range("C11").Value = 406

# Row 16 - new entry "scottie"
$ws.Range("A16").Value = "scottie"
$ws.Range("B16").Value = 26
$ws.Range("C16").Value = 406

# Row 17 - new entry "john"
$ws.Range("A17").Value = "john"
$ws.Range("B17").Value = 28
$ws.Range("C17").Value = 470

# Row 18 - new entry "billybob" (no previous payment value)
$ws.Range("A18").Value = "billybob"
$ws.Range("B18").Value = 12

# Row 19 - new entry "bobbert"
$ws.Range("A19").Value = "bobbert"
$ws.Range("B19").Value = 12

# Update selection to match the new view
$ws.Range("B3").Select() | Out-Null
